$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.48654696685555
$ws.Range("C2").Value = 8.240977258690252
$ws.Range("D2").Value = 9.235474447145522
$ws.Range("E2").Value = 13.51009371198987
$ws.Range("F2").Value = 31.73970379576939
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 9.944062962671387
$ws.Range("M2").Value = 16.71181436925978
$ws.Range("O2").Value = 23.80790824098542
$ws.Range("B3").Value = 13.91101951493269
$ws.Range("C3").Value = 7.730390033144063
$ws.Range("D3").Value = 9.216359259248568
$ws.Range("E3").Value = 13.53106918420188
$ws.Range("F3").Value = 31.83984451700831
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 9.974609070204284
$ws.Range("M3").Value = 16.52247938479737
$ws.Range("O3").Value = 23.91556842161789
$ws.Range("B4").Value = 13.54625247433391
$ws.Range("C4").Value = 7.3979282640722
$ws.Range("D4").Value = 9.205663590728118
$ws.Range("E4").Value = 13.54613840086088
$ws.Range("F4").Value = 31.91096345801212
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 9.994593257764627
$ws.Range("M4").Value = 16.40714083840787
$ws.Range("O4").Value = 23.98855442043322
$ws.Range("B5").Value = 13.39495098139617
$ws.Range("C5").Value = 7.257699424194038
$ws.Range("D5").Value = 9.201569485162034
$ws.Range("E5").Value = 13.55282988678308
$ws.Range("F5").Value = 31.94235665550525
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 10.00304635976567
$ws.Range("M5").Value = 16.36041412054749
$ws.Range("O5").Value = 24.0200202566911
$ws.Range("B6").Value = 13.36967364434066
$ws.Range("C6").Value = 7.234128484294719
$ws.Range("D6").Value = 9.200905709831193
$ws.Range("E6").Value = 13.55397425717054
$ws.Range("F6").Value = 31.94771481285344
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 10.0044686888588
$ws.Range("M6").Value = 16.35267307713022
$ws.Range("O6").Value = 24.02534901580601
$ws.Range("B7").Value = 13.54422242994892
$ws.Range("C7").Value = 7.396056282142149
$ws.Range("D7").Value = 9.205607301969517
$ws.Range("E7").Value = 13.54622641516125
$ws.Range("F7").Value = 31.91137708725532
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 9.994706006037141
$ws.Range("M7").Value = 16.40650949444543
$ws.Range("O7").Value = 23.98897181109142
$ws.Range("B8").Value = 14.29058675869829
$ws.Range("C8").Value = 8.068871208454933
$ws.Range("D8").Value = 9.228669086768724
$ws.Range("E8").Value = 13.51687145352213
$ws.Range("F8").Value = 31.77222708533209
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 9.954340474956453
$ws.Range("M8").Value = 16.64636901303066
$ws.Range("O8").Value = 23.84359681464541
$ws.Range("B9").Value = 15.65587662853101
$ws.Range("C9").Value = 9.23731482250213
$ws.Range("D9").Value = 9.282031425263956
$ws.Range("E9").Value = 13.47668615368232
$ws.Range("F9").Value = 31.57621933887826
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 9.884915177762753
$ws.Range("M9").Value = 17.12200196913667
$ws.Range("O9").Value = 23.61343776194103
$ws.Range("B10").Value = 16.58992355463608
$ws.Range("C10").Value = 10.00344169797898
$ws.Range("D10").Value = 9.326025543749008
$ws.Range("E10").Value = 13.45775577746594
$ws.Range("F10").Value = 31.47962799472348
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 9.839815672502402
$ws.Range("M10").Value = 17.47200939203137
$ws.Range("O10").Value = 23.47824134432572
$ws.Range("B11").Value = 16.99833506246378
$ws.Range("C11").Value = 10.3319376910695
$ws.Range("D11").Value = 9.347039575943098
$ws.Range("E11").Value = 13.45144205011673
$ws.Range("F11").Value = 31.44608545964117
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 9.820576209307747
$ws.Range("M11").Value = 17.63080181295783
$ws.Range("O11").Value = 23.42418512782362
$ws.Range("B12").Value = 17.15051473763639
$ws.Range("C12").Value = 10.45345570203355
$ws.Range("D12").Value = 9.355137320339033
$ws.Range("E12").Value = 13.44938124640444
$ws.Range("F12").Value = 31.43488605235828
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 9.813473917002383
$ws.Range("M12").Value = 17.69082512370643
$ws.Range("O12").Value = 23.40479295468537
$ws.Range("B13").Value = 17.11785180156799
$ws.Range("C13").Value = 10.42741248761477
$ws.Range("D13").Value = 9.353387152424098
$ws.Range("E13").Value = 13.44981040427002
$ws.Range("F13").Value = 31.43723111780273
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 9.814995376027396
$ws.Range("M13").Value = 17.67790355863959
$ws.Range("O13").Value = 23.40892136226844
$ws.Range("B14").Value = 17.01090516702183
$ws.Range("C14").Value = 10.3419926749508
$ws.Range("D14").Value = 9.347702994627905
$ws.Range("E14").Value = 13.4512658939084
$ws.Range("F14").Value = 31.44513392826705
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 9.819988227965263
$ws.Range("M14").Value = 17.63574239236167
$ws.Range("O14").Value = 23.42256808251965
$ws.Range("B15").Value = 16.94507182906415
$ws.Range("C15").Value = 10.28929610242238
$ws.Range("D15").Value = 9.344239426847469
$ws.Range("E15").Value = 13.45220039560207
$ws.Range("F15").Value = 31.45017048893873
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 9.82307035332823
$ws.Range("M15").Value = 17.60990204365365
$ws.Range("O15").Value = 23.43106765410067
$ws.Range("B16").Value = 16.56289168773316
$ws.Range("C16").Value = 9.981571036366542
$ws.Range("D16").Value = 9.324672051313543
$ws.Range("E16").Value = 13.45821460426536
$ws.Range("F16").Value = 31.48202994357966
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 9.841098678834131
$ws.Range("M16").Value = 17.46161933630906
$ws.Range("O16").Value = 23.48192447715953
$ws.Range("B17").Value = 16.32413275469779
$ws.Range("C17").Value = 9.787665353699989
$ws.Range("D17").Value = 9.312921836318026
$ws.Range("E17").Value = 13.46249243060586
$ws.Range("F17").Value = 31.50424317059765
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 9.852485214181931
$ws.Range("M17").Value = 17.37050981741296
$ws.Range("O17").Value = 23.51503554513942
$ws.Range("B18").Value = 16.18525797003018
$ws.Range("C18").Value = 9.674251354537846
$ws.Range("D18").Value = 9.306257864786639
$ws.Range("E18").Value = 13.46516923046758
$ws.Range("F18").Value = 31.51799789787255
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 9.859154617685743
$ws.Range("M18").Value = 17.31806849659514
$ws.Range("O18").Value = 23.53478053109831
$ws.Range("B19").Value = 16.13797513107638
$ws.Range("C19").Value = 9.635527607141988
$ws.Range("D19").Value = 9.304017890508396
$ws.Range("E19").Value = 13.46611270886386
$ws.Range("F19").Value = 31.52282279364491
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 9.861433411280972
$ws.Range("M19").Value = 17.30030775003256
$ws.Range("O19").Value = 23.5415859139313
$ws.Range("B20").Value = 16.34970996828626
$ws.Range("C20").Value = 9.808501980708391
$ws.Range("D20").Value = 9.314162919081378
$ws.Range("E20").Value = 13.46201466390617
$ws.Range("F20").Value = 31.50177723201598
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 9.851260663070734
$ws.Range("M20").Value = 17.38021280540435
$ws.Range("O20").Value = 23.51143827692577
$ws.Range("B21").Value = 17.04238597387433
$ws.Range("C21").Value = 10.367160576747
$ws.Range("D21").Value = 9.349368795446864
$ws.Range("E21").Value = 13.45082942672123
$ws.Range("F21").Value = 31.44277185223561
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 9.818516734550977
$ws.Range("M21").Value = 17.64812944320008
$ws.Range("O21").Value = 23.41853040472595
$ws.Range("B22").Value = 17.4806139484905
$ws.Range("C22").Value = 10.71551353248097
$ws.Range("D22").Value = 9.373193151523461
$ws.Range("E22").Value = 13.44544287496433
$ws.Range("F22").Value = 31.41296844032298
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 9.79818482281369
$ws.Range("M22").Value = 17.82257854905935
$ws.Range("O22").Value = 23.36409479198397
$ws.Range("B23").Value = 17.24807897779548
$ws.Range("C23").Value = 10.53112347937823
$ws.Range("D23").Value = 9.36040431787344
$ws.Range("E23").Value = 13.44814190969016
$ws.Range("F23").Value = 31.42807136292396
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 9.808938705060337
$ws.Range("M23").Value = 17.72954631502137
$ws.Range("O23").Value = 23.39257072701384
$ws.Range("B24").Value = 16.33815151383504
$ws.Range("C24").Value = 9.799087769089654
$ws.Range("D24").Value = 9.313601540646335
$ws.Range("E24").Value = 13.46222998496318
$ws.Range("F24").Value = 31.50288901811583
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 9.851813899049487
$ws.Range("M24").Value = 17.37582627174576
$ws.Range("O24").Value = 23.51306239390532
$ws.Range("B25").Value = 15.29807495653768
$ws.Range("C25").Value = 8.937427633310714
$ws.Range("D25").Value = 9.266741412232674
$ws.Range("E25").Value = 13.48569647954042
$ws.Range("F25").Value = 31.62095528387196
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 9.902657390386672
$ws.Range("M25").Value = 16.99304129737314
$ws.Range("O25").Value = 24.0200202566911
